# Apply the "Recall notice elements:" header row + switch recall markup from
# <span> to <p>, driven off a new $H$45 reference cell, and reorder/relabel
# the notice field rows 46-56.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- sheetView: scroll position / selection as captured in the saved file ---
$ws.Application.ActiveWindow.ScrollRow = 21
$ws.Application.ActiveWindow.ScrollColumn = 9
$ws.Range("J30").Select()

# --- New row 45: label + element-tag cells (bold, like K44/L44) ---
$ws.Range("G45").Value = "Recall notice elements:"
$ws.Range("H45").Value = "p"
$ws.Range("G45:H45").Font.Bold = $true

# --- Field order/labels for rows 46-56 ---
$fields = @(
    @{Row=46; Field="Manufacturer";        Label="Manufacturer"},
    @{Row=47; Field="ModelYear";           Label="Model Year"},
    @{Row=48; Field="Make";                Label="Make"},
    @{Row=49; Field="Model";               Label="Model"},
    @{Row=50; Field="NHTSACampaignNumber"; Label="NHTSA Campaign Number"},
    @{Row=51; Field="ReportReceivedDate";  Label="Report Received Date"},
    @{Row=52; Field="Component";           Label="Component"},
    @{Row=53; Field="Summary";             Label="Summary"},
    @{Row=54; Field="Conequence";          Label="Consequence"},
    @{Row=55; Field="Remedy";              Label="Remedy"},
    @{Row=56; Field="Notes";               Label="Notes"}
)

foreach ($item in $fields) {
    $r = $item.Row

    $ws.Range("F$r").Value = $item.Field
    $ws.Range("G$r").Value = $item.Label

    $ws.Range("E$r").Formula = '="<"&$H$45&" class=""recall-heading"" >"'
    $ws.Range("H$r").Formula = '="</"&$H$45&">"'
    $ws.Range("I$r").Formula = '="<"&$H$45&" class=""recall-data"" >"'
    $ws.Range("J$r").Formula = '="{{ notice."&F' + $r + '&" ? "&"notice."&F' + $r + '&" : ""n/a"" }}"'
    $ws.Range("K$r").Formula = '="</"&$H$45&">"'
    $ws.Range("L$r").Formula = '="<div class=""" & $L$44&""" >"&E' + $r + '&G' + $r + '&H' + $r + '&I' + $r + '&J' + $r + '&K' + $r + '&"</div>"'
}
